$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100: Nigeria
$ws.Cells.Item(100,1).Value = "Nigeria"
$ws.Cells.Item(100,2).Value = 442
$ws.Cells.Item(100,3).Value = 35
$ws.Cells.Item(100,4).Value = 152
$ws.Cells.Item(100,5).Value = 277
$ws.Cells.Item(100,6).Value = 2
$ws.Cells.Item(100,7).Value = 1
$ws.Cells.Item(100,8).Value = 13

# Row 101: Bolivia
$ws.Cells.Item(101,1).Value = "Bolivia"
$ws.Cells.Item(101,2).Value = 441
$ws.Cells.Item(101,3).Value = 44
$ws.Cells.Item(101,4).Value = 14
$ws.Cells.Item(101,5).Value = 398
$ws.Cells.Item(101,6).Value = 3
$ws.Cells.Item(101,7).Value = 1
$ws.Cells.Item(101,8).Value = 29

# Row 102: Guinea
$ws.Cells.Item(102,1).Value = "Guinea"
$ws.Cells.Item(102,2).Value = 438
$ws.Cells.Item(102,3).Value = 34
$ws.Cells.Item(102,4).Value = 49
$ws.Cells.Item(102,5).Value = 388
$ws.Cells.Item(102,6).Value = 0
$ws.Cells.Item(102,7).Value = 0
$ws.Cells.Item(102,8).Value = 1

# Row 103: Honduras
$ws.Cells.Item(103,1).Value = "Honduras"
$ws.Cells.Item(103,2).Value = 426
$ws.Cells.Item(103,3).Value = 7
$ws.Cells.Item(103,4).Value = 9
$ws.Cells.Item(103,5).Value = 382
$ws.Cells.Item(103,6).Value = 10
$ws.Cells.Item(103,7).Value = 4
$ws.Cells.Item(103,8).Value = 35

# Row 104: San Marino
$ws.Cells.Item(104,1).Value = "San Marino"
$ws.Cells.Item(104,2).Value = 426
$ws.Cells.Item(104,3).Value = 33
$ws.Cells.Item(104,4).Value = 55
$ws.Cells.Item(104,5).Value = 333
$ws.Cells.Item(104,6).Value = 15
$ws.Cells.Item(104,7).Value = 2
$ws.Cells.Item(104,8).Value = 38

# Row 105: Malta
$ws.Cells.Item(105,1).Value = "Malta"
$ws.Cells.Item(105,2).Value = 412
$ws.Cells.Item(105,3).Value = 13
$ws.Cells.Item(105,4).Value = 82
$ws.Cells.Item(105,5).Value = 327
$ws.Cells.Item(105,6).Value = 4
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 3

# Row 108: Reunion
$ws.Cells.Item(108,1).Value = "Reunion"
$ws.Cells.Item(108,2).Value = 394
$ws.Cells.Item(108,3).Value = 3
$ws.Cells.Item(108,4).Value = 237
$ws.Cells.Item(108,5).Value = 157
$ws.Cells.Item(108,6).Value = 4
$ws.Cells.Item(108,7).Value = 0
$ws.Cells.Item(108,8).Value = 0

# Row 177: Nueva Caledonia
$ws.Cells.Item(177,1).Value = "Nueva Caledonia"
$ws.Cells.Item(177,2).Value = 18
$ws.Cells.Item(177,3).Value = 0
$ws.Cells.Item(177,4).Value = 1
$ws.Cells.Item(177,5).Value = 17
$ws.Cells.Item(177,6).Value = 1
$ws.Cells.Item(177,7).Value = 0
$ws.Cells.Item(177,8).Value = 0

# Row 178: Timor Oriental
$ws.Cells.Item(178,1).Value = "Timor Oriental"
$ws.Cells.Item(178,2).Value = 18
$ws.Cells.Item(178,3).Value = 10
$ws.Cells.Item(178,4).Value = 1
$ws.Cells.Item(178,5).Value = 17
$ws.Cells.Item(178,6).Value = 0
$ws.Cells.Item(178,7).Value = 0
$ws.Cells.Item(178,8).Value = 0

# Row 180: Islas Virgenes de los Estados Unidos
$ws.Cells.Item(180,1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(180,2).Value = 17
$ws.Cells.Item(180,3).Value = 0
$ws.Cells.Item(180,4).Value = 0
$ws.Cells.Item(180,5).Value = 17
$ws.Cells.Item(180,6).Value = 0
$ws.Cells.Item(180,7).Value = 0
$ws.Cells.Item(180,8).Value = 0

# Row 181: Fiyi
$ws.Cells.Item(181,1).Value = "Fiyi"
$ws.Cells.Item(181,2).Value = 17
$ws.Cells.Item(181,3).Value = 1
$ws.Cells.Item(181,4).Value = 0
$ws.Cells.Item(181,5).Value = 17
$ws.Cells.Item(181,6).Value = 0
$ws.Cells.Item(181,7).Value = 0
$ws.Cells.Item(181,8).Value = 0

# Row 182: Nepal
$ws.Cells.Item(182,1).Value = "Nepal"
$ws.Cells.Item(182,2).Value = 16
$ws.Cells.Item(182,3).Value = 0
$ws.Cells.Item(182,4).Value = 2
$ws.Cells.Item(182,5).Value = 14
$ws.Cells.Item(182,6).Value = 0
$ws.Cells.Item(182,7).Value = 0
$ws.Cells.Item(182,8).Value = 0

# Row 183: Malaui
$ws.Cells.Item(183,1).Value = "Malaui"
$ws.Cells.Item(183,2).Value = 16
$ws.Cells.Item(183,3).Value = 0
$ws.Cells.Item(183,4).Value = 0
$ws.Cells.Item(183,5).Value = 14
$ws.Cells.Item(183,6).Value = 1
$ws.Cells.Item(183,7).Value = 0
$ws.Cells.Item(183,8).Value = 2

# Row 197: Islas Malvinas
$ws.Cells.Item(197,1).Value = "Islas Malvinas"
$ws.Cells.Item(197,2).Value = 11
$ws.Cells.Item(197,3).Value = 0
$ws.Cells.Item(197,4).Value = 1
$ws.Cells.Item(197,5).Value = 10
$ws.Cells.Item(197,6).Value = 0
$ws.Cells.Item(197,7).Value = 0
$ws.Cells.Item(197,8).Value = 0

# Row 198: Montserrat
$ws.Cells.Item(198,1).Value = "Montserrat"
$ws.Cells.Item(198,2).Value = 11
$ws.Cells.Item(198,3).Value = 0
$ws.Cells.Item(198,4).Value = 1
$ws.Cells.Item(198,5).Value = 10
$ws.Cells.Item(198,6).Value = 1
$ws.Cells.Item(198,7).Value = 0
$ws.Cells.Item(198,8).Value = 0

# Row 215: Yemen
$ws.Cells.Item(215,1).Value = "Yemen"
$ws.Cells.Item(215,2).Value = 1
$ws.Cells.Item(215,3).Value = 0
$ws.Cells.Item(215,4).Value = 0
$ws.Cells.Item(215,5).Value = 1
$ws.Cells.Item(215,6).Value = 0
$ws.Cells.Item(215,7).Value = 0
$ws.Cells.Item(215,8).Value = 0

# Row 216: San Pedro y Miquelon
$ws.Cells.Item(216,1).Value = "San Pedro y Miquelon"
$ws.Cells.Item(216,2).Value = 1
$ws.Cells.Item(216,3).Value = 0
$ws.Cells.Item(216,4).Value = 0
$ws.Cells.Item(216,5).Value = 1
$ws.Cells.Item(216,6).Value = 0
$ws.Cells.Item(216,7).Value = 0
$ws.Cells.Item(216,8).Value = 0
